# Applies the "break out stock.yaml completed" update to the weekly
# KPITTECH.NS history sheet:
#   - Q64: 1 -> 0
#   - O269: 0 -> 2
#   - R271, R272: blank -> 0
#   - appends 26 new weekly rows (273-298), extending the used range to R298

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- targeted fixes on existing rows ---------------------------------
$ws.Cells.Item(64, 17).Value = 0    # Q64
$ws.Cells.Item(269, 15).Value = 2   # O269
$ws.Cells.Item(271, 18).Value = 0   # R271
$ws.Cells.Item(272, 18).Value = 0   # R272

# --- new weekly rows 273-298 ------------------------------------------
# columns: row, A Datetime, B Open, C High, D Low, E Close, G Volume,
#          H Year, I Month, J Day, K Hour, L Minute, M Second, N Week,
#          O isPivot, P two_line_structure, Q detect_structure
# (F "Adj Close" and R "backup" are left blank for every new row, matching
# the source data)
$newRows = @(
    @(273,45474,1645.556945321405,1769.322749406045,1631.594643967163,1689.388549804688,7245721,2024,7,1,0,0,0,27,0,0,0),
    @(274,45481,1689.388587774959,1923.506467530388,1649.296884539043,1860.526611328125,9587488,2024,7,8,0,0,0,28,1,0,0),
    @(275,45488,1877.231457059516,1903.859511536272,1805.126095275355,1814.799926757812,6510024,2024,7,15,0,0,0,29,0,1,1),
    @(276,45495,1795.152983832635,1880.422750564685,1700.957366238724,1784.132690429688,6500989,2024,7,22,0,0,0,30,0,0,0),
    @(277,45502,1795.152992966442,1904.856786981058,1769.422418037282,1819.038452148438,5159465,2024,7,29,0,0,0,31,0,0,0),
    @(278,45509,1764.785123748103,1790.116630947117,1680.462745591653,1753.515502929688,3544267,2024,8,5,0,0,0,32,0,0,0),
    @(279,45516,1743,1855.949951171875,1716,1815.400024414062,3402361,2024,8,12,0,0,0,33,0,0,0),
    @(280,45523,1831.900024414062,1899,1801,1826.5,2909822,2024,8,19,0,0,0,34,0,0,2),
    @(281,45530,1845,1921,1807,1815.150024414062,4126338,2024,8,26,0,0,0,35,1,0,0),
    @(282,45537,1815.150024414062,1824.300048828125,1715.449951171875,1730.300048828125,3357562,2024,9,2,0,0,0,36,0,0,0),
    @(283,45544,1718.25,1859.949951171875,1707.150024414062,1826.050048828125,3103712,2024,9,9,0,0,0,37,0,0,0),
    @(284,45551,1832,1832,1626.099975585938,1662,4075198,2024,9,16,0,0,0,38,0,0,0),
    @(285,45558,1680,1713.699951171875,1630,1671,5640959,2024,9,23,0,0,0,39,0,0,0),
    @(286,45565,1637.800048828125,1713.300048828125,1625,1691.25,4643500,2024,9,30,0,0,0,40,0,0,0),
    @(287,45572,1699,1797,1641.800048828125,1789.050048828125,3609958,2024,10,7,0,0,0,41,0,2,2),
    @(288,45579,1790,1822.900024414062,1750.050048828125,1786.349975585938,2531418,2024,10,14,0,0,0,42,0,0,0),
    @(289,45586,1790,1797.599975585938,1338.449951171875,1373.599975585938,18055673,2024,10,21,0,0,0,43,0,0,0),
    @(290,45593,1385.949951171875,1429.900024414062,1309.050048828125,1424.25,7924376,2024,10,28,0,0,0,44,0,0,0),
    @(291,45600,1430.050048828125,1520,1378,1422.349975585938,6845019,2024,11,4,0,0,0,45,0,0,0),
    @(292,45607,1422.349975585938,1447.449951171875,1339,1356.949951171875,3346714,2024,11,11,0,0,0,46,0,0,0),
    @(293,45614,1350.25,1356.400024414062,1283.25,1308.449951171875,4656848,2024,11,18,0,0,0,47,2,0,0),
    @(294,45621,1308.449951171875,1424.5,1308.449951171875,1368.800048828125,6260192,2024,11,25,0,0,0,48,0,0,0),
    @(295,45628,1361,1503.75,1340,1481.099975585938,12355971,2024,12,2,0,0,0,49,0,0,0),
    @(296,45635,1482,1563.349975585938,1480.550048828125,1533.099975585938,6076152,2024,12,9,0,0,0,50,0,0,0),
    @(297,45642,1535,1555,1438.900024414062,1451.25,4704310,2024,12,16,0,0,0,51,0,0,0),
    @(298,45649,1467.400024414062,1474.849975585938,1408.050048828125,1460.599975585938,2694935,2024,12,23,0,0,0,52,0,0,0)
)

# Column A carries the same custom datetime display format as the rest of
# the sheet's Datetime column.
$ws.Range("A273:A298").NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value  = $r[1]   # A Datetime
    $ws.Cells.Item($row, 2).Value  = $r[2]   # B Open
    $ws.Cells.Item($row, 3).Value  = $r[3]   # C High
    $ws.Cells.Item($row, 4).Value  = $r[4]   # D Low
    $ws.Cells.Item($row, 5).Value  = $r[5]   # E Close
    $ws.Cells.Item($row, 7).Value  = $r[6]   # G Volume
    $ws.Cells.Item($row, 8).Value  = $r[7]   # H Year
    $ws.Cells.Item($row, 9).Value  = $r[8]   # I Month
    $ws.Cells.Item($row, 10).Value = $r[9]   # J Day
    $ws.Cells.Item($row, 11).Value = $r[10]  # K Hour
    $ws.Cells.Item($row, 12).Value = $r[11]  # L Minute
    $ws.Cells.Item($row, 13).Value = $r[12]  # M Second
    $ws.Cells.Item($row, 14).Value = $r[13]  # N Week
    $ws.Cells.Item($row, 15).Value = $r[14]  # O isPivot
    $ws.Cells.Item($row, 16).Value = $r[15]  # P two_line_structure
    $ws.Cells.Item($row, 17).Value = $r[16]  # Q detect_structure
}
